$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Update header row (row 1) values for columns B-E
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Update row 2 values for columns B-E
$ws.Range("B2").Value = 12.327871893356299
$ws.Range("C2").Value = 39.249740724287321
$ws.Range("D2").Value = 51.556526058029348
$ws.Range("E2").Value = 46.228756715520007

# Update row 3 values for columns B-E
$ws.Range("B3").Value = 28.211836691470481
$ws.Range("C3").Value = 68.194223218324865
$ws.Range("D3").Value = 68.282344595357159
$ws.Range("E3").Value = 43.044488056068573

# Update the selection on the sheet to match B1:E3
$ws.Range("B1:E3").Select()
